# "Validation step for the structuring"
#
# The "Fasciolo S1_SS1_B1_001" mapping row (Serie/Sottoserie/Busta/Fascicolo
# rdf:type rows) on the Immagini(temp) sheet was cleared out while validating
# the RDF structuring - the "Immagini..." record-resource row right below it
# (and everything further down the sheet) is left untouched.

$wb = $excel.ActiveWorkbook

$immagini = $wb.Worksheets.Item("Immagini(temp)")
$immagini.Activate()
$immagini.Range("A2:G5").ClearContents()

# Leave the cursor where the author ended up validating each sheet.
$sottoserie = $wb.Worksheets.Item("Sottoserie")
$sottoserie.Activate()
$sottoserie.Range("E12").Select()

$documento = $wb.Worksheets.Item("Documento")
$documento.Activate()
$documento.Range("F12").Select()

$immagini.Activate()
$immagini.Range("B19").Select()
